$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 17:22"

# Row 4
$ws.Range("B4").Value = 506008
$ws.Range("C4").Value = 3132
$ws.Range("D4").Value = 28210
$ws.Range("E4").Value = 458915
$ws.Range("G4").Value = 136
$ws.Range("H4").Value = 18883

# Row 16
$ws.Range("B16").Value = 22559
$ws.Range("C16").Value = 411
$ws.Range("E16").Value = 15977

# Row 17
$ws.Range("B17").Value = 20022
$ws.Range("C17").Value = 233
$ws.Range("E17").Value = 18774

# Row 19
$ws.Range("B19").Value = 13789
$ws.Range("C19").Value = 229
$ws.Range("E19").Value = 6848

# Row 27
$ws.Range("B27").Value = 6927
$ws.Range("C27").Value = 426
$ws.Range("D27").Value = 1864
$ws.Range("E27").Value = 4990
$ws.Range("F27").Value = 383
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = 73

# Row 54
$ws.Range("B54").Value = 2081
$ws.Range("C54").Value = 70
$ws.Range("E54").Value = 1719
$ws.Range("G54").Value = 1
$ws.Range("H54").Value = 93

# Row 85
$ws.Range("A85").Value = "Republica de Chipre"
$ws.Range("B85").Value = 616
$ws.Range("C85").Value = 21
$ws.Range("D85").Value = 58
$ws.Range("E85").Value = 548
$ws.Range("F85").Value = 11
$ws.Range("H85").Value = 10

# Row 86
$ws.Range("A86").Value = "Principado de Andorra"
$ws.Range("B86").Value = 601
$ws.Range("D86").Value = 71
$ws.Range("E86").Value = 504
$ws.Range("F86").Value = 17
$ws.Range("H86").Value = 26

# Row 99
$ws.Range("A99").Value = "Reunion"
$ws.Range("B99").Value = 388
$ws.Range("C99").Value = 6
$ws.Range("D99").Value = 40
$ws.Range("E99").Value = 348
$ws.Range("F99").Value = 3
$ws.Range("H99").Value = 0

# Row 100
$ws.Range("A100").Value = "Taiwan"
$ws.Range("B100").Value = 385
$ws.Range("C100").Value = 3
$ws.Range("D100").Value = 99
$ws.Range("E100").Value = 280
$ws.Range("F100").Value = 0
$ws.Range("H100").Value = 6
